# Mise à jour de l'application
# Add a new daily attendance column (CC) for 2025-11-19, mirroring the
# formatting of the preceding day column (CB) and filling in each
# player's attendance mark for the new day.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

function Set-DayCell($row, $value) {
    $dst = $ws.Range("CC" + $row)
    # Write the value first so dependent COUNTA/COUNTIF formulas on the
    # row pick it up during recalculation, then clone the neighbouring
    # day cell's number format / alignment onto it.
    if ($null -ne $value) {
        $dst.Value2 = $value
    }
    $ws.Range("CB" + $row).Copy()
    $dst.PasteSpecial($xlPasteFormats)
}

# Header: new date column (19 Nov 2025 -> serial 45980), same format as CB1
Set-DayCell 1 45980

# Attendance marks per player row.
# "P" = Présent, "B" = Blessure, $null = leave the cell blank (style-only,
# matching row 21 whose attendance tracking already stopped at CB21).
$marks = [ordered]@{
    2  = "P"
    3  = "P"
    4  = "P"
    5  = "B"
    6  = "B"
    7  = "P"
    8  = "P"
    9  = "P"
    10 = "P"
    11 = "P"
    13 = "B"
    14 = "P"
    15 = "P"
    16 = "P"
    17 = "P"
    18 = "P"
    19 = "P"
    20 = "P"
    21 = $null
    22 = "P"
    23 = "P"
    24 = "P"
    25 = "P"
    26 = "P"
    27 = "P"
    28 = "P"
    29 = "P"
}

foreach ($row in $marks.Keys) {
    Set-DayCell $row $marks[$row]
}

# Row 12's records stop at column AX (player left mid-season), so it gets
# no CC cell - nothing to do for that row.

# Move the active selection to CC24, matching the author's final cursor spot.
[void]$ws.Range("CC24").Select()
